# Regenerate the s_vals data (espino_paolo, 2021) to filter save games.
# Updates the numeric statistic columns (TB, d2S, K, IP, sum) for every
# data row while leaving the date column (A) and Win column (F) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2021-06-23)
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 5.582307763322248

# Row 3 (2021-05-14)
$ws.Range("B3").Value = 0.2881169905109251
$ws.Range("C3").Value = 109.9114832445916
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 124.7845412238902

# Row 4 (2021-05-11)
$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 19.48425592650926

# Row 5 (2021-05-05)
$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 109.9114832445916
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 127.7687514715588

# Row 6 (2021-05-04)
$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 0.1496068669990043
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 5.582307763322248

# Row 7 (2021-04-24)
$ws.Range("B7").Value = 3.272327238179451
$ws.Range("C7").Value = 1.626987699542094
$ws.Range("D7").Value = 0.1496068669990043
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 5.582307763322248
